$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New card rows (10-13) -------------------------------------------------
# Shared-string insertion order matters for byte-identical sharedStrings.xml,
# so values are written in the same order the original author must have
# used: I10, B10, B11, B12, I11, I12, B13, I13 (new text values), with the
# purely-numeric cells filled in alongside each row.

# Row 10 - Fonio Millet
$ws.Range("I10").Value = "100857.png"
$ws.Range("B10").Value = "Fonio Millet"
$ws.Range("C10").Value = 2021.01
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 48
$ws.Range("F10").Value = 6768
$ws.Range("G10").Value = 7.03
$ws.Range("H10").Value = 100857

# Row 11 - Hazelnut
$ws.Range("B11").Value = "Hazelnut"
$ws.Range("C11").Value = 2021.03
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 35
$ws.Range("F11").Value = 353.56
$ws.Range("G11").Value = 1.21
$ws.Range("H11").Value = 100877

# Row 12 - Yellowhorn Tree
$ws.Range("B12").Value = "Yellowhorn Tree"
$ws.Range("C12").Value = 2019.03
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 35
$ws.Range("F12").Value = 485.86
$ws.Range("G12").Value = 2.24
$ws.Range("H12").Value = 100606

# Back-fill the image-file-name cells for rows 11 & 12
$ws.Range("I11").Value = "100877.png"
$ws.Range("I12").Value = "100606.jpg"

# Row 13 - Pink Ipe(circumflex) Tree
$ws.Range("B13").Value = "Pink Ipê Tree"
$ws.Range("C13").Value = 2017.11
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 140.1
$ws.Range("G13").Value = 0.2
$ws.Range("H13").Value = 100379
$ws.Range("I13").Value = "100379.jpg"

# --- K2:K13 shared CONCATENATE formula --------------------------------------
$ws.Range("K2:K13").Formula = "=CONCATENATE(""const cards"",A2,"" = new theCards('"",B2,""', "",C2,"", "",D2,"", "",E2,"", "",F2,"", "",G2,"", 'https://doi.org/10.5524/"",H2,""' ,'./images/"",I2,""');"")"

# --- New L/M data columns (file-size histogram) -----------------------------
$ws.Range("L12").Value = 7776
$ws.Range("M12").Value = 2

$ws.Range("L13").Value = 991966016
$ws.Range("M13").Value = 4

$ws.Range("L14").Value = 116854
$ws.Range("M14").Value = 3

$ws.Range("L15").Value = 33971884
$ws.Range("M15").Value = 1

$ws.Range("L16").Value = 11801015
$ws.Range("M16").Value = 1

$ws.Range("L17").Value = 3786774
$ws.Range("M17").Value = 14

$ws.Range("L18").Value = 22780139
$ws.Range("M18").Value = 1

$ws.Range("L19").Value = 256182
$ws.Range("M19").Value = 6

$ws.Range("L20").Value = 589
$ws.Range("M20").Value = 1

$ws.Range("L21").Value = 6580680
$ws.Range("M21").Value = 2

$ws.Range("L22").Value = 1333274198
$ws.Range("M22").Value = 4

# --- Totals / unit-conversion formulas --------------------------------------
$ws.Range("L36").Formula = "=SUM(L12:L35)"
$ws.Range("L37").Formula = "=L36/1024"
$ws.Range("L38").Formula = "=L37/1024"
$ws.Range("L39").Formula = "=L38/1024"

# --- View / selection / print setup -----------------------------------------
$ws.Range("I14").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
